$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("control_panel")
$ws.Activate()

# Clear the scenario numbers that used to let GeneXpert be toggled via scenario values
$ws.Range("B10").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()

# Update default_smoothness value
$ws.Range("B17").Value = 0.01

# Move the active selection to B11 to match the saved view state
$ws.Range("B11").Select()
